# Add a new slide (slide 13) "Pour aller plus loin (quartz)" after the
# existing last slide, using the same "Title and Content" layout (layout 2)
# as the other content slides in this deck (e.g. the previous slide,
# "Spring scheduler").

$p = $ppt.ActivePresentation

$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Match the shape naming convention used by the rest of the (French
# locale) deck.
$s.Shapes.Item(1).Name = "Titre 1"
$s.Shapes.Item(2).Name = "Espace réservé du contenu 2"

# --- Title -----------------------------------------------------------
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Pour aller plus loin (quartz)"

# --- Body content ------------------------------------------------------
$tf = $s.Shapes.Item(2).TextFrame

$body = "La facilité de planifier des exécutions à travers les annotations a une contrepartie : "
$body += "`rDéclarations statiques des planifications"
$body += "`rPlanification non distribuée !"
$body += "`r"
$body += "`rL’intégration de Quartz est la solution ! Mais …"
$body += "`rUtilisation d’une base de données"
$body += "`rArchitecture et concepts techniques + complexe"
$body += "`rPlus de tuyauterie "
$body += "`r"

$tf.TextRange.Text = $body

$paras = $tf.TextRange.Paragraphs()

# Paragraph 2: "Déclarations statiques des planifications" - centered, red
$para2 = $tf.TextRange.Paragraphs(2, 1)
$para2.ParagraphFormat.Alignment = 2
$para2.Font.Color.RGB = 255

# Paragraph 3: "Planification non distribuée !" - centered, red
$para3 = $tf.TextRange.Paragraphs(3, 1)
$para3.ParagraphFormat.Alignment = 2
$para3.Font.Color.RGB = 255

# Paragraphs 6,7,8: bulleted list items
foreach ($i in 6, 7, 8) {
    $para = $tf.TextRange.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Visible = -1
    $para.ParagraphFormat.Bullet.Character = 8226
    $para.ParagraphFormat.Bullet.Font.Name = "Arial"
}

Write-Host "Added slide with $($s.Shapes.Count) shapes; total slides: $($p.Slides.Count)"
